$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = [double]"6624.091939615049"
$ws.Range("C2").Value = [double]"1.152690925201376e-20"
$ws.Range("B3").Value = [double]"228.3432100179251"
$ws.Range("C3").Value = [double]"0.7253843033282803"
$ws.Range("B4").Value = [double]"849.7573214361937"
$ws.Range("C4").Value = [double]"0.1457770827067722"
$ws.Range("B5").Value = [double]"1535.922120336277"
$ws.Range("C5").Value = [double]"0.008493861009430104"
$ws.Range("B6").Value = [double]"-62.61052269387119"
$ws.Range("C6").Value = [double]"0.5908718352635725"
$ws.Range("B7").Value = [double]"-258.0216149905932"
$ws.Range("C7").Value = [double]"0.009858628649475147"
$ws.Range("B8").Value = [double]"78.02413245656567"
$ws.Range("C8").Value = [double]"0.006990887960258856"
$ws.Range("B9").Value = [double]"-1277.030284844331"
$ws.Range("C9").Value = [double]"1.289458636219465e-64"
$ws.Range("B10").Value = [double]"-32.84530079177129"
$ws.Range("C10").Value = [double]"4.826909747181773e-24"
$ws.Range("B11").Value = [double]"239.0369259702444"
$ws.Range("C11").Value = [double]"6.30923317330064e-30"
$ws.Range("B12").Value = [double]"445.4485496797773"
$ws.Range("C12").Value = [double]"1.818890819261763e-280"
$ws.Range("B13").Value = [double]"-0.02413144128049513"
$ws.Range("C13").Value = [double]"0.04694911845140084"
$ws.Range("B14").Value = [double]"3.255668937909595e-05"
$ws.Range("C14").Value = [double]"0.02730461039220844"
$ws.Range("B15").Value = [double]"-25.83877854850373"
$ws.Range("C15").Value = [double]"1.596093540208e-23"
$ws.Range("B16").Value = [double]"7.24351491082813"
$ws.Range("C16").Value = [double]"0.004313573001920735"
$ws.Range("B17").Value = [double]"-1966.248626537106"
$ws.Range("C17").Value = [double]"1.4318288446911e-07"
$ws.Range("B18").Value = [double]"-149.9071073948375"
$ws.Range("C18").Value = [double]"0.5810298995962568"
$ws.Name = "summ34555815"

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = [double]"6770.416538600455"
$ws.Range("C2").Value = [double]"7.481096794173274e-21"
$ws.Range("B3").Value = [double]"36.71413177178229"
$ws.Range("C3").Value = [double]"0.9555010747403475"
$ws.Range("B4").Value = [double]"462.9651397558928"
$ws.Range("C4").Value = [double]"0.4392297329630611"
$ws.Range("B5").Value = [double]"1219.922269252028"
$ws.Range("C5").Value = [double]"0.04135134155520814"
$ws.Range("B6").Value = [double]"-94.6889818414727"
$ws.Range("C6").Value = [double]"0.4183051311744916"
$ws.Range("B7").Value = [double]"-227.751984452231"
$ws.Range("C7").Value = [double]"0.02316714409058028"
$ws.Range("B8").Value = [double]"113.5538912029445"
$ws.Range("C8").Value = [double]"9.097129418943146e-05"
$ws.Range("B9").Value = [double]"-1340.182653460834"
$ws.Range("C9").Value = [double]"6.49053520592413e-71"
$ws.Range("B10").Value = [double]"-29.47804287150083"
$ws.Range("C10").Value = [double]"1.173917301664931e-19"
$ws.Range("B11").Value = [double]"197.7280554508511"
$ws.Range("C11").Value = [double]"4.532913515899602e-21"
$ws.Range("B12").Value = [double]"439.7447124841459"
$ws.Range("C12").Value = [double]"1.302857052492772e-273"
$ws.Range("B13").Value = [double]"-0.0225813169988327"
$ws.Range("C13").Value = [double]"0.06134256784186077"
$ws.Range("B14").Value = [double]"2.878566165154093e-05"
$ws.Range("C14").Value = [double]"0.05097137887155125"
$ws.Range("B15").Value = [double]"-24.1258948483392"
$ws.Range("C15").Value = [double]"7.206748656493802e-21"
$ws.Range("B16").Value = [double]"8.794377387350583"
$ws.Range("C16").Value = [double]"0.0004885401152967951"
$ws.Range("B17").Value = [double]"-2231.34871761781"
$ws.Range("C17").Value = [double]"2.624884254325082e-09"
$ws.Range("B18").Value = [double]"-393.8029110088944"
$ws.Range("C18").Value = [double]"0.1459828770894746"
$ws.Name = "summ35753127"

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = [double]"7461.884322494781"
$ws.Range("C2").Value = [double]"5.608152917964229e-26"
$ws.Range("B3").Value = [double]"-283.7418021681182"
$ws.Range("C3").Value = [double]"0.6592626827162009"
$ws.Range("B4").Value = [double]"186.3786852737967"
$ws.Range("C4").Value = [double]"0.7480920009085562"
$ws.Range("B5").Value = [double]"893.6644498136789"
$ws.Range("C5").Value = [double]"0.1231638496105257"
$ws.Range("B6").Value = [double]"-114.8958036545955"
$ws.Range("C6").Value = [double]"0.3249537074210901"
$ws.Range("B7").Value = [double]"-257.2481792334563"
$ws.Range("C7").Value = [double]"0.0102516230227415"
$ws.Range("B8").Value = [double]"93.15787536609291"
$ws.Range("C8").Value = [double]"0.00134714471489296"
$ws.Range("B9").Value = [double]"-1360.382298968127"
$ws.Range("C9").Value = [double]"9.080878894725828e-73"
$ws.Range("B10").Value = [double]"-32.24890916056569"
$ws.Range("C10").Value = [double]"5.244779507303432e-23"
$ws.Range("B11").Value = [double]"217.6541517995254"
$ws.Range("C11").Value = [double]"5.367097120166904e-25"
$ws.Range("B12").Value = [double]"445.859520305143"
$ws.Range("C12").Value = [double]"3.485212817436992e-278"
$ws.Range("B13").Value = [double]"-0.02008937285535889"
$ws.Range("C13").Value = [double]"0.09621111740061811"
$ws.Range("B14").Value = [double]"2.3174340070814e-05"
$ws.Range("C14").Value = [double]"0.1167883344416816"
$ws.Range("B15").Value = [double]"-24.1258948483392"
$ws.Range("C15").Value = [double]"7.206748656493802e-21"
$ws.Range("B16").Value = [double]"5.065632204777085"
$ws.Range("C16").Value = [double]"0.04651687030465003"
$ws.Range("B17").Value = [double]"-1667.959216028098"
$ws.Range("C17").Value = [double]"9.194723969321303e-06"
$ws.Range("B18").Value = [double]"-428.6262390009616"
$ws.Range("C18").Value = [double]"0.1145259578556049"
$ws.Name = "summ36811193"

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = [double]"6823.916890347897"
$ws.Range("C2").Value = [double]"1.286827525840077e-21"
$ws.Range("B3").Value = [double]"91.62280008119797"
$ws.Range("C3").Value = [double]"0.88824617749574"
$ws.Range("B4").Value = [double]"671.6669506372211"
$ws.Range("C4").Value = [double]"0.2537434100175722"
$ws.Range("B5").Value = [double]"1337.527887343113"
$ws.Range("C5").Value = [double]"0.0229206158050227"
$ws.Range("B6").Value = [double]"42.96210738965419"
$ws.Range("C6").Value = [double]"0.7131509122532754"
$ws.Range("B7").Value = [double]"-128.3664969536162"
$ws.Range("C7").Value = [double]"0.2004451600718109"
$ws.Range("B8").Value = [double]"104.3508366290975"
$ws.Range("C8").Value = [double]"0.0003203755862749113"
$ws.Range("B9").Value = [double]"-1258.607936368335"
$ws.Range("C9").Value = [double]"1.526301797157229e-62"
$ws.Range("B10").Value = [double]"-30.60434505871893"
$ws.Range("C10").Value = [double]"5.795755020322535e-21"
$ws.Range("B11").Value = [double]"215.179312857918"
$ws.Range("C11").Value = [double]"1.584129843721568e-24"
$ws.Range("B12").Value = [double]"435.9591805685413"
$ws.Range("C12").Value = [double]"3.585343340254555e-267"
$ws.Range("B13").Value = [double]"-0.02494793216347008"
$ws.Range("C13").Value = [double]"0.03991983012623504"
$ws.Range("B14").Value = [double]"3.57912271018716e-05"
$ws.Range("C14").Value = [double]"0.01554169046054736"
$ws.Range("B15").Value = [double]"-25.92947680689074"
$ws.Range("C15").Value = [double]"1.752725652642588e-23"
$ws.Range("B16").Value = [double]"6.258740397548511"
$ws.Range("C16").Value = [double]"0.01398315529817734"
$ws.Range("B17").Value = [double]"-2204.984718017239"
$ws.Range("C17").Value = [double]"4.237487337528569e-09"
$ws.Range("B18").Value = [double]"-279.5707451670112"
$ws.Range("C18").Value = [double]"0.3011945415957252"
$ws.Name = "summ37848199"

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = [double]"6772.266371410433"
$ws.Range("C2").Value = [double]"3.091495950978651e-22"
$ws.Range("B3").Value = [double]"148.6227433579546"
$ws.Range("C3").Value = [double]"0.8142633474790536"
$ws.Range("B4").Value = [double]"576.4355839728762"
$ws.Range("C4").Value = [double]"0.3101396459399253"
$ws.Range("B5").Value = [double]"1318.586939378448"
$ws.Range("C5").Value = [double]"0.02010821168176978"
$ws.Range("B6").Value = [double]"-36.9606154081423"
$ws.Range("C6").Value = [double]"0.7518734240366013"
$ws.Range("B7").Value = [double]"-223.8808695505404"
$ws.Range("C7").Value = [double]"0.02555819603430334"
$ws.Range("B8").Value = [double]"99.40429453026802"
$ws.Range("C8").Value = [double]"0.0005817998790250989"
$ws.Range("B9").Value = [double]"-1247.681187692295"
$ws.Range("C9").Value = [double]"1.322528067265362e-61"
$ws.Range("B10").Value = [double]"-32.41558636150799"
$ws.Range("C10").Value = [double]"2.327395389574531e-23"
$ws.Range("B11").Value = [double]"230.4473649741129"
$ws.Range("C11").Value = [double]"8.399444023428662e-28"
$ws.Range("B12").Value = [double]"442.0104853269651"
$ws.Range("C12").Value = [double]"5.417882395395106e-274"
$ws.Range("B13").Value = [double]"-0.0298520183727372"
$ws.Range("C13").Value = [double]"0.01447272834016086"
$ws.Range("B14").Value = [double]"3.880068859382592e-05"
$ws.Range("C14").Value = [double]"0.009162854342617512"
$ws.Range("B15").Value = [double]"-24.90550945871123"
$ws.Range("C15").Value = [double]"8.384415631488711e-22"
$ws.Range("B16").Value = [double]"6.759928558929783"
$ws.Range("C16").Value = [double]"0.007490433889891894"
$ws.Range("B17").Value = [double]"-1751.417888830631"
$ws.Range("C17").Value = [double]"3.047381424635422e-06"
$ws.Range("B18").Value = [double]"-282.5811309649564"
$ws.Range("C18").Value = [double]"0.2991243610071354"
$ws.Name = "summ38871872"

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = [double]"7296.663555991818"
$ws.Range("C2").Value = [double]"3.757937689641432e-23"
$ws.Range("B3").Value = [double]"-626.3237708660671"
$ws.Range("C3").Value = [double]"0.3557746386386572"
$ws.Range("B4").Value = [double]"-129.9816084846104"
$ws.Range("C4").Value = [double]"0.8330405860996495"
$ws.Range("B5").Value = [double]"693.2418940406366"
$ws.Range("C5").Value = [double]"0.2603692758361625"
$ws.Range("B6").Value = [double]"-33.76555384022868"
$ws.Range("C6").Value = [double]"0.7718973038997154"
$ws.Range("B7").Value = [double]"-207.6063360360661"
$ws.Range("C7").Value = [double]"0.03728785028398279"
$ws.Range("B8").Value = [double]"80.05253548812401"
$ws.Range("C8").Value = [double]"0.005620277546041102"
$ws.Range("B9").Value = [double]"-1265.680615212637"
$ws.Range("C9").Value = [double]"2.377903781552092e-63"
$ws.Range("B10").Value = [double]"-30.86929268565315"
$ws.Range("C10").Value = [double]"3.653241024219768e-21"
$ws.Range("B11").Value = [double]"215.3671611559883"
$ws.Range("C11").Value = [double]"1.934973261091925e-24"
$ws.Range("B12").Value = [double]"464.5138668965926"
$ws.Range("C12").Value = [double]"3.344671169035491e-300"
$ws.Range("B13").Value = [double]"-0.01773280665650193"
$ws.Range("C13").Value = [double]"0.1474258826896417"
$ws.Range("B14").Value = [double]"2.292331543134412e-05"
$ws.Range("C14").Value = [double]"0.1287518731481276"
$ws.Range("B15").Value = [double]"-23.64216014670562"
$ws.Range("C15").Value = [double]"4.005349368841608e-20"
$ws.Range("B16").Value = [double]"7.700907990028298"
$ws.Range("C16").Value = [double]"0.001931809311898648"
$ws.Range("B17").Value = [double]"-1844.655651692567"
$ws.Range("C17").Value = [double]"9.355277287036075e-07"
$ws.Range("B18").Value = [double]"-462.5388360868985"
$ws.Range("C18").Value = [double]"0.08910411961319634"
$ws.Name = "summ39881069"

$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = [double]"6223.50974287706"
$ws.Range("C2").Value = [double]"5.049042998306949e-19"
$ws.Range("B3").Value = [double]"206.6713868423355"
$ws.Range("C3").Value = [double]"0.7464864878610826"
$ws.Range("B4").Value = [double]"782.751000691729"
$ws.Range("C4").Value = [double]"0.1702770731126856"
$ws.Range("B5").Value = [double]"1483.214139199416"
$ws.Range("C5").Value = [double]"0.0092832380023826"
$ws.Range("B6").Value = [double]"-77.21334820449215"
$ws.Range("C6").Value = [double]"0.5091113329083763"
$ws.Range("B7").Value = [double]"-223.4985329474478"
$ws.Range("C7").Value = [double]"0.02617329955271449"
$ws.Range("B8").Value = [double]"91.23018685231938"
$ws.Range("C8").Value = [double]"0.001688678772838237"
$ws.Range("B9").Value = [double]"-1201.9205400017"
$ws.Range("C9").Value = [double]"2.795417698386581e-57"
$ws.Range("B10").Value = [double]"-32.4363413093256"
$ws.Range("C10").Value = [double]"2.739529208681029e-23"
$ws.Range("B11").Value = [double]"216.5011678113882"
$ws.Range("C11").Value = [double]"7.860593787168775e-25"
$ws.Range("B12").Value = [double]"463.5862043768341"
$ws.Range("C12").Value = [double]"3.303842775133711e-301"
$ws.Range("B13").Value = [double]"-0.0312384323633382"
$ws.Range("C13").Value = [double]"0.01037391156737841"
$ws.Range("B14").Value = [double]"3.381731478430233e-05"
$ws.Range("C14").Value = [double]"0.02369478619780072"
$ws.Range("B15").Value = [double]"-23.71272487543796"
$ws.Range("C15").Value = [double]"6.197890969005883e-20"
$ws.Range("B16").Value = [double]"8.846712799532639"
$ws.Range("C16").Value = [double]"0.0004384083763430969"
$ws.Range("B17").Value = [double]"-1480.054684185312"
$ws.Range("C17").Value = [double]"7.840180706595151e-05"
$ws.Range("B18").Value = [double]"-221.0081882479477"
$ws.Range("C18").Value = [double]"0.4161061458262869"
$ws.Name = "summ40903302"

$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = [double]"6778.323971204394"
$ws.Range("C2").Value = [double]"4.682774967396987e-22"
$ws.Range("B3").Value = [double]"-211.806706859307"
$ws.Range("C3").Value = [double]"0.7423065920013899"
$ws.Range("B4").Value = [double]"367.8986827526811"
$ws.Range("C4").Value = [double]"0.5237665969545897"
$ws.Range("B5").Value = [double]"1103.828304497855"
$ws.Range("C5").Value = [double]"0.05551665643462327"
$ws.Range("B6").Value = [double]"-55.02837299913745"
$ws.Range("C6").Value = [double]"0.6356089477782173"
$ws.Range("B7").Value = [double]"-117.7111544732977"
$ws.Range("C7").Value = [double]"0.2369485278690177"
$ws.Range("B8").Value = [double]"109.1030312158527"
$ws.Range("C8").Value = [double]"0.0001755285764330724"
$ws.Range("B9").Value = [double]"-1178.883165316118"
$ws.Range("C9").Value = [double]"2.083497091335132e-55"
$ws.Range("B10").Value = [double]"-30.88065426810451"
$ws.Range("C10").Value = [double]"2.787650763915135e-21"
$ws.Range("B11").Value = [double]"202.1182063705942"
$ws.Range("C11").Value = [double]"9.561140384348961e-22"
$ws.Range("B12").Value = [double]"446.7741323709686"
$ws.Range("C12").Value = [double]"3.638130397458268e-280"
$ws.Range("B13").Value = [double]"-0.02639447795500077"
$ws.Range("C13").Value = [double]"0.02987725243512843"
$ws.Range("B14").Value = [double]"3.727041043448277e-05"
$ws.Range("C14").Value = [double]"0.01243445143050265"
$ws.Range("B15").Value = [double]"-24.46205756312416"
$ws.Range("C15").Value = [double]"1.7373709030365e-21"
$ws.Range("B16").Value = [double]"8.39990067573178"
$ws.Range("C16").Value = [double]"0.0007112519274289875"
$ws.Range("B17").Value = [double]"-2032.062406501771"
$ws.Range("C17").Value = [double]"5.416096613617508e-08"
$ws.Range("B18").Value = [double]"-421.123064619183"
$ws.Range("C18").Value = [double]"0.119683113082327"
$ws.Name = "summ41944001"

$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = [double]"6469.865614423885"
$ws.Range("C2").Value = [double]"6.848341149141642e-18"
$ws.Range("B3").Value = [double]"-476.9148289161361"
$ws.Range("C3").Value = [double]"0.488441144701511"
$ws.Range("B4").Value = [double]"-99.81275351369584"
$ws.Range("C4").Value = [double]"0.8740177008336181"
$ws.Range("B5").Value = [double]"710.8134993487113"
$ws.Range("C5").Value = [double]"0.2583434574887769"
$ws.Range("B6").Value = [double]"-120.1862164420497"
$ws.Range("C6").Value = [double]"0.3033930943800523"
$ws.Range("B7").Value = [double]"-211.4804593348385"
$ws.Range("C7").Value = [double]"0.03479434624262429"
$ws.Range("B8").Value = [double]"108.335424797817"
$ws.Range("C8").Value = [double]"0.0001958494675141304"
$ws.Range("B9").Value = [double]"-1253.634580223938"
$ws.Range("C9").Value = [double]"4.215017714181354e-62"
$ws.Range("B10").Value = [double]"-31.88122448690157"
$ws.Range("C10").Value = [double]"1.359533371565941e-22"
$ws.Range("B11").Value = [double]"261.7703611497037"
$ws.Range("C11").Value = [double]"7.797171244032146e-35"
$ws.Range("B12").Value = [double]"455.2287478640305"
$ws.Range("C12").Value = [double]"2.634591681065067e-287"
$ws.Range("B13").Value = [double]"-0.02019380558666431"
$ws.Range("C13").Value = [double]"0.09553924873788543"
$ws.Range("B14").Value = [double]"2.616403421535114e-05"
$ws.Range("C14").Value = [double]"0.0768786233147175"
$ws.Range("B15").Value = [double]"-20.80546747949741"
$ws.Range("C15").Value = [double]"1.138340654470172e-15"
$ws.Range("B16").Value = [double]"11.82637353805136"
$ws.Range("C16").Value = [double]"4.221567452558846e-06"
$ws.Range("B17").Value = [double]"-1480.272157185571"
$ws.Range("C17").Value = [double]"7.406634751459634e-05"
$ws.Range("B18").Value = [double]"-310.2107332385683"
$ws.Range("C18").Value = [double]"0.2549510942898031"
$ws.Name = "summ43063316"

Write-Host "Done updating all sheets"